# Applies the Fenrir_Profits market-data refresh described by the commit diff.
# All touched cells are plain numeric literals (no formulas) on the per-server leve tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 770.0909
$ws.Range("I6").Value = 763.44446
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 2290.33338
$ws.Range("L6").Value = 2400
$ws.Range("M6").Value = -2178.33338
$ws.Range("N6").Value = -2624
$ws.Range("H18").Value = 920.5
$ws.Range("I18").Value = 920.5
$ws.Range("K18").Value = 920.5
$ws.Range("M18").Value = -636.5
$ws.Range("H19").Value = 713.2857
$ws.Range("I19").Value = 699.2
$ws.Range("J19").Value = 721.1111
$ws.Range("K19").Value = 699.2
$ws.Range("L19").Value = 721.1111
$ws.Range("M19").Value = -524.2
$ws.Range("N19").Value = -1071.1111
$ws.Range("H33").Value = 69.11539
$ws.Range("I33").Value = 64.91304
$ws.Range("J33").Value = 101.333336
$ws.Range("K33").Value = 64.91304
$ws.Range("L33").Value = 101.333336
$ws.Range("M33").Value = 164.08696
$ws.Range("N33").Value = -559.333336
$ws.Range("H41").Value = 413.92307
$ws.Range("I41").Value = 366.55554
$ws.Range("J41").Value = 520.5
$ws.Range("K41").Value = 366.55554
$ws.Range("L41").Value = 520.5
$ws.Range("M41").Value = 73.44445999999999
$ws.Range("N41").Value = -1400.5
$ws.Range("H53").Value = 218.75
$ws.Range("I53").Value = 85.5
$ws.Range("J53").Value = 352
$ws.Range("K53").Value = 85.5
$ws.Range("L53").Value = 352
$ws.Range("M53").Value = 551.5
$ws.Range("N53").Value = -1626
$ws.Range("H74").Value = 6821.231
$ws.Range("I74").Value = 7067.6
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 7067.6
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -6131.6
$ws.Range("N74").Value = -7872
$ws.Range("H77").Value = 6821.231
$ws.Range("I77").Value = 7067.6
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 35338
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -30658
$ws.Range("N77").Value = -39360
$ws.Range("H98").Value = 17547042
$ws.Range("I98").Value = 18521796
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 18521796
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -18520298
$ws.Range("N98").Value = -4496
$ws.Range("H122").Value = 17547042
$ws.Range("I122").Value = 18521796
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 55565388
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -55562938
$ws.Range("N122").Value = -9400
$ws.Range("H123").Value = 41653.332
$ws.Range("J123").Value = 41653.332
$ws.Range("L123").Value = 41653.332
$ws.Range("N123").Value = -51453.332
$ws.Range("H124").Value = 54780
$ws.Range("J124").Value = 54780
$ws.Range("L124").Value = 54780
$ws.Range("N124").Value = -64600
$ws.Range("H125").Value = 918.6
$ws.Range("J125").Value = 990.25
$ws.Range("L125").Value = 8912.25
$ws.Range("N125").Value = -13832.25
$ws.Range("H126").Value = 40701.43
$ws.Range("J126").Value = 40701.43
$ws.Range("L126").Value = 40701.43
$ws.Range("N126").Value = -50581.43
$ws.Range("H127").Value = 2849.5386
$ws.Range("I127").Value = 1940
$ws.Range("J127").Value = 3122.4
$ws.Range("K127").Value = 5820
$ws.Range("L127").Value = 9367.200000000001
$ws.Range("M127").Value = -860
$ws.Range("N127").Value = -19287.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 35800708
$ws.Range("I80").Value = 183682.67
$ws.Range("J80").Value = 62513476
$ws.Range("K80").Value = 183682.67
$ws.Range("L80").Value = 62513476
$ws.Range("M80").Value = -182684.67
$ws.Range("N80").Value = -62515472
$ws.Range("H83").Value = 35800708
$ws.Range("I83").Value = 183682.67
$ws.Range("J83").Value = 62513476
$ws.Range("K83").Value = 918413.3500000001
$ws.Range("L83").Value = 312567380
$ws.Range("M83").Value = -913421.3500000001
$ws.Range("N83").Value = -312577364
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H98").Value = 32998.332
$ws.Range("J98").Value = 32998.332
$ws.Range("L98").Value = 32998.332
$ws.Range("N98").Value = -38988.332
$ws.Range("H105").Value = 1746.6666
$ws.Range("I105").Value = 500
$ws.Range("J105").Value = 2993.3333
$ws.Range("K105").Value = 500
$ws.Range("L105").Value = 2993.3333
$ws.Range("M105").Value = 1247
$ws.Range("N105").Value = -6487.3333
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3878.261
$ws.Range("J80").Value = 4104.7617
$ws.Range("L80").Value = 12314.2851
$ws.Range("N80").Value = -14186.2851
$ws.Range("H83").Value = 3878.261
$ws.Range("J83").Value = 4104.7617
$ws.Range("L83").Value = 36942.8553
$ws.Range("N83").Value = -46302.8553
$ws.Range("H122").Value = 742.7857
$ws.Range("I122").Value = 427.18182
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 3844.63638
$ws.Range("L122").Value = 17100
$ws.Range("M122").Value = -1394.63638
$ws.Range("N122").Value = -22000
$ws.Range("H123").Value = 3913
$ws.Range("I123").Value = 1618.75
$ws.Range("J123").Value = 4747.273
$ws.Range("K123").Value = 4856.25
$ws.Range("L123").Value = 14241.819
$ws.Range("M123").Value = -2406.25
$ws.Range("N123").Value = -19141.819
$ws.Range("H124").Value = 3223.2856
$ws.Range("I124").Value = 3030
$ws.Range("J124").Value = 3255.5
$ws.Range("K124").Value = 9090
$ws.Range("L124").Value = 9766.5
$ws.Range("M124").Value = -4180
$ws.Range("N124").Value = -19586.5
$ws.Range("H125").Value = 6000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H126").Value = 1566.25
$ws.Range("H129").Value = 848.2917
$ws.Range("I129").Value = 369.875
$ws.Range("K129").Value = 1109.625
$ws.Range("M129").Value = 3890.375
$ws.Range("H130").Value = 1444.4445
$ws.Range("I130").Value = 1000
$ws.Range("K130").Value = 3000
$ws.Range("M130").Value = 2020
$ws.Range("H131").Value = 20441248
$ws.Range("I131").Value = 50000396
$ws.Range("J131").Value = 2526611.2
$ws.Range("K131").Value = 150001188
$ws.Range("L131").Value = 7579833.600000001
$ws.Range("M131").Value = -149996148
$ws.Range("N131").Value = -7589913.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26959.158
$ws.Range("I70").Value = 37283.266
$ws.Range("J70").Value = 4836.0713
$ws.Range("K70").Value = 37283.266
$ws.Range("L70").Value = 4836.0713
$ws.Range("M70").Value = -37013.266
$ws.Range("N70").Value = -5376.0713
$ws.Range("H73").Value = 26959.158
$ws.Range("I73").Value = 37283.266
$ws.Range("J73").Value = 4836.0713
$ws.Range("K73").Value = 37283.266
$ws.Range("L73").Value = 4836.0713
$ws.Range("M73").Value = -36347.266
$ws.Range("N73").Value = -6708.0713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 34133
$ws.Range("I16").Value = 1199
$ws.Range("K16").Value = 1199
$ws.Range("M16").Value = -1029
$ws.Range("H55").Value = 4786.12
$ws.Range("I55").Value = 1001.7222
$ws.Range("K55").Value = 1001.7222
$ws.Range("M55").Value = -828.7222

